# Applies the OOXML diff:
#  1. Removes the "sic" comment (w:id=2, author Celine Camps) that wrapped
#     the "m" in "...trasse premierem" and replaces its comment markup with
#     inline <corr><exp>ent</exp></corr> correction markup (regenerated
#     from a new download, per the commit message).
#  2. Splits " le cane avecq un " into " le ca" + "v" + "e avecq un ",
#     i.e. turns "cane" into "cave" with the new "v" in slightly
#     different (un-colored) formatting.

$d = $word.ActiveDocument

# ---- Step 0: capture formatting templates BEFORE any mutation happens ----
# (FormattedText captures both text + full run formatting; copying from an
#  existing, identically-styled run is the only reliable way in this host
#  to get every rFonts/sz/szCs attribute written out, matching what Word
#  itself would serialize.)

# "<corr>" / "</corr>" use the same styling as the existing "<del>" / "</del>"
# tags elsewhere in the document (Courier New, color a91111, sz18/szCs18).
$tmplCorrOpen  = $d.Range(745, 750).FormattedText   # "<del>"  -> becomes "<corr>"
$tmplCorrClose = $d.Range(756, 762).FormattedText   # "</del>" -> becomes "</corr>" (clean, no strike)

# "<exp>" / "ent" / "</exp>" already exist verbatim elsewhere in the document
# (color a9a9a9 Courier New sz14/szCs14 for the tags, plain black for "ent").
$tmplExpOpen   = $d.Range(366, 371).FormattedText   # "<exp>"
$tmplEnt       = $d.Range(371, 374).FormattedText   # "ent"
$tmplExpClose  = $d.Range(374, 380).FormattedText   # "</exp>"

# Minimal (rtl-only, no explicit color/font) run style, used for the new "v".
$tmplPlain     = $d.Range(144, 145).FormattedText   # "T"

# ---- Step 1: delete the obsolete comment (Celine Camps, "sic") ----
# This removes the commentRangeStart/commentRangeEnd/commentReference markup
# around the "m" as well as the <w:comment> definition in comments.xml.
for ($i = 1; $i -le $d.Comments.Count; $i++) {
    $cmt = $d.Comments.Item($i)
    if ($cmt.Author -eq "Celine Camps") {
        $cmt.Delete()
        break
    }
}

# ---- Step 2: insert "<corr><exp>ent</exp></corr>" right after the "m" ----
$rngM = $d.Content
$rngM.Find.ClearFormatting()
$rngM.Find.Execute("trasse premierem", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insPos = $rngM.End

$d.Range($insPos, $insPos).FormattedText = $tmplCorrOpen
$d.Range($insPos, $insPos + 5).Text = "<corr>"
$insPos = $insPos + 6

$d.Range($insPos, $insPos).FormattedText = $tmplExpOpen
$d.Range($insPos, $insPos + 5).Text = "<exp>"
$insPos = $insPos + 5

$d.Range($insPos, $insPos).FormattedText = $tmplEnt
$d.Range($insPos, $insPos + 3).Text = "ent"
$insPos = $insPos + 3

$d.Range($insPos, $insPos).FormattedText = $tmplExpClose
$d.Range($insPos, $insPos + 6).Text = "</exp>"
$insPos = $insPos + 6

$d.Range($insPos, $insPos).FormattedText = $tmplCorrClose
$d.Range($insPos, $insPos + 6).Text = "</corr>"
$insPos = $insPos + 7

# ---- Step 3: split " le cane avecq un " -> " le ca" + "v" + "e avecq un " ----
$rngFind = $d.Content
$rngFind.Find.ClearFormatting()
$rngFind.Find.Execute("cane avecq", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$nPos = $rngFind.Start + 2

$d.Range($nPos, $nPos + 1).FormattedText = $tmplPlain
$d.Range($nPos, $nPos + 1).Text = "v"
